# Fix a big error in the diagram:
#  1. Slide 1 title is corrected from "AddressBook – Level 4" to "Doer List".
#  2. The stray blank slide 9 at the end of the deck is removed.

$p = $ppt.ActivePresentation

# --- 1. Fix the title text on the first slide ------------------------------
# The title currently reads two runs: "AddressBook" (flagged as a misspelling)
# followed by " - Level 4". Drop the first ("AddressBook") run entirely and
# replace the remaining text with the corrected title "Doer List", keeping
# the surviving run's (non-flagged) formatting.
$s1 = $p.Slides.Item(1)
$title = $s1.Shapes.Item(1)
$tr = $title.TextFrame.TextRange
$firstRun = $tr.Characters(1, 11)
$firstRun.Text = ""
$tr.Text = "Doer List"

# --- 2. Remove the trailing empty slide ------------------------------------
$p.Slides.Item($p.Slides.Count).Delete()
